$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '47.424.89'
$ws.Range("E2").Value = '  +2.69%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.508.65'
$ws.Range("E3").Value = '  +2.10%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.10'
$ws.Range("E5").Value = '  +0.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.54'
$ws.Range("E6").Value = '  +3.84%  '

$ws.Range("E7").Value = '  +1.55%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.542'
$ws.Range("E9").Value = '  +0.97%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.05'
$ws.Range("E10").Value = '  +8.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0819'
$ws.Range("E11").Value = '  +1.37%  '

$ws.Range("E12").Value = '  +0.97%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.67'
$ws.Range("E13").Value = '  +1.64%  '

$ws.Range("E14").Value = '  +2.00%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.901.59'
$ws.Range("E15").Value = '  +1.94%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.515.40'
$ws.Range("E16").Value = '  +2.84%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.864'
$ws.Range("E17").Value = '  +2.34%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '47.385.81'
$ws.Range("E18").Value = '  +2.80%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.03'
$ws.Range("E19").Value = '  +3.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.73'
$ws.Range("E20").Value = '  +4.23%  '

$ws.Range("E21").Value = '  +0.87%  '

$ws.Range("B22").Value = 'ImmutableX'
$ws.Range("C22").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.63'
$ws.Range("E22").Value = '  +11.52%  '

$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.04'
$ws.Range("E23").Value = '  -1.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '249.89'
$ws.Range("E24").Value = '  +0.67%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.61'
$ws.Range("E25").Value = '  +3.76%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.22'
$ws.Range("E26").Value = '  +0.64%  '

$ws.Range("E27").Value = '  -0.05%  '

$ws.Range("E28").Value = '  +4.87%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.07'
$ws.Range("E29").Value = '  +3.80%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.03'
$ws.Range("E30").Value = '  +6.01%  '

$ws.Range("E31").Value = '  +4.47%  '

$ws.Range("E32").Value = '  +1.89%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.96'
$ws.Range("E33").Value = '  -2.10%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.46'
$ws.Range("E34").Value = '  +2.67%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0795'
$ws.Range("E35").Value = '  +4.05%  '

$ws.Range("E36").Value = '  +0.16%  '

$ws.Range("E37").Value = '  +5.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.76'
$ws.Range("E38").Value = '  +4.23%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.01'
$ws.Range("E39").Value = '  +2.62%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '122.75'
$ws.Range("E41").Value = '  -4.42%  '

$ws.Range("E42").Value = '  -1.26%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.49'
$ws.Range("E43").Value = '  +2.64%  '

$ws.Range("E44").Value = '  +2.19%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.996.31'
$ws.Range("E45").Value = '  +1.84%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.09'
$ws.Range("E46").Value = '  +3.54%  '

$ws.Range("E47").Value = '  -1.44%  '

$ws.Range("E48").Value = '  -2.79%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.08'
$ws.Range("E49").Value = '  -1.45%  '

$ws.Range("E50").Value = '  +8.11%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.64'
$ws.Range("E51").Value = '  +1.25%  '
